$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "326.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.16%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.39%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.568"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-6.62%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08091"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.52%"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.71%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.336"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.66%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.905"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.61%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9496"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.99%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1162"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.07%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1897"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.54%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1008"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.87%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04174"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "4.92%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1063"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.07%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001270"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.13%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006073"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.13%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.616"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.41%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.539"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-6.19%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1374"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.06%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04272"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.42%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001237"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.61%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004617"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.42%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001234"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.43%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004000"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.03%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02673"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.01%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05559"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.55%"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "24.65%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007699"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.92%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1395"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.99%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002076"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.40%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008707"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.83%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007113"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.66%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.06%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003488"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-5.00%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002276"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.30%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.06%"
